$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 4
